$d = $word.ActiveDocument

# --- 1. Locate the paragraph that contains "left-top" (inside <ab><margin>left-top</margin>) ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*left-top*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate paragraph containing 'left-top'"
}

# --- 2. Replace that whole paragraph with two new paragraphs:
#      - same text but "left-middle" instead of "left-top", with the pBdr/shd
#        paragraph-border cruft dropped and the run-formatting folded into pPr/rPr
#      - a brand new paragraph containing "<render>tall</render>"
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml =
  '<w:p ' + $xmlNs + '>' +
    '<w:pPr>' +
      '<w:widowControl w:val="0"/>' +
      '<w:contextualSpacing w:val="0"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' +
        '<w:color w:val="7f6000"/>' +
        '<w:sz w:val="18"/>' +
        '<w:szCs w:val="18"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' +
        '<w:color w:val="7f6000"/>' +
        '<w:sz w:val="18"/>' +
        '<w:szCs w:val="18"/>' +
        '<w:rtl w:val="0"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">&lt;ab&gt;&lt;margin&gt;</w:t>' +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rtl w:val="0"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">left-middle</w:t>' +
    '</w:r>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' +
        '<w:color w:val="7f6000"/>' +
        '<w:sz w:val="18"/>' +
        '<w:szCs w:val="18"/>' +
        '<w:rtl w:val="0"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">&lt;/margin&gt;</w:t>' +
    '</w:r>' +
  '</w:p>' +
  '<w:p ' + $xmlNs + '>' +
    '<w:pPr>' +
      '<w:widowControl w:val="0"/>' +
      '<w:contextualSpacing w:val="0"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' +
        '<w:color w:val="7f6000"/>' +
        '<w:sz w:val="18"/>' +
        '<w:szCs w:val="18"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' +
        '<w:color w:val="7f6000"/>' +
        '<w:sz w:val="18"/>' +
        '<w:szCs w:val="18"/>' +
        '<w:rtl w:val="0"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">&lt;render&gt;tall&lt;/render&gt;</w:t>' +
    '</w:r>' +
  '</w:p>'

$null = $target.Range.InsertXML($newXml)

# --- 3. sectPr / pgMar: add an explicit footer distance of 720 twips (36pt) ---
$d.PageSetup.FooterDistance = 36
